$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.08982198144684
$ws.Range("C2").Value = 8.442805355930762
$ws.Range("D2").Value = 7.561073078221929
$ws.Range("F2").Value = 39.19551131286678
$ws.Range("G2").Value = 45.73263836846137
$ws.Range("H2").Value = 18.49891040264304
$ws.Range("I2").Value = 28.6133229083567
$ws.Range("J2").Value = 10.2849378823748
$ws.Range("L2").Value = 11.9214355413239
$ws.Range("M2").Value = 17.40694084634886
$ws.Range("B3").Value = 17.66818391846313
$ws.Range("C3").Value = 8.019438413419179
$ws.Range("D3").Value = 7.560786092553042
$ws.Range("F3").Value = 39.28936685984774
$ws.Range("G3").Value = 45.7939770487658
$ws.Range("H3").Value = 18.5593245895113
$ws.Range("I3").Value = 28.74338401402257
$ws.Range("J3").Value = 10.30159963985198
$ws.Range("L3").Value = 11.92822296155308
$ws.Range("M3").Value = 17.32174758916222
$ws.Range("B4").Value = 17.4069612433295
$ws.Range("C4").Value = 7.74620958989696
$ws.Range("D4").Value = 7.560880406984824
$ws.Range("F4").Value = 39.3572124746735
$ws.Range("G4").Value = 45.84705306802417
$ws.Range("H4").Value = 18.60021177369272
$ws.Range("I4").Value = 28.82957640436146
$ws.Range("J4").Value = 10.31236273761893
$ws.Range("L4").Value = 11.93377002426
$ws.Range("M4").Value = 17.27165983887944
$ws.Range("B5").Value = 17.30007595397914
$ws.Range("C5").Value = 7.631597293467757
$ws.Range("D5").Value = 7.560987126908891
$ws.Range("F5").Value = 39.38742054029444
$ws.Range("G5").Value = 45.8725412455187
$ws.Range("H5").Value = 18.61782527715954
$ws.Range("I5").Value = 28.86628995665529
$ws.Range("J5").Value = 10.31688316097339
$ws.Range("L5").Value = 11.93637809460386
$ws.Range("M5").Value = 17.25182156696818
$ws.Range("B6").Value = 17.28230603094715
$ws.Range("C6").Value = 7.612370862127124
$ws.Range("D6").Value = 7.561008979591937
$ws.Range("F6").Value = 39.39259093499302
$ws.Range("G6").Value = 45.87700609534228
$ws.Range("H6").Value = 18.62080740218996
$ws.Range("I6").Value = 28.87248211525576
$ws.Range("J6").Value = 10.31764190199284
$ws.Range("L6").Value = 11.93683217594979
$ws.Range("M6").Value = 17.24856244390495
$ws.Range("B7").Value = 17.40552130637118
$ws.Range("C7").Value = 7.74467702339808
$ws.Range("D7").Value = 7.560881569402356
$ws.Range("F7").Value = 39.35760951655574
$ws.Range("G7").Value = 45.84738120833593
$ws.Range("H7").Value = 18.60044546516876
$ws.Range("I7").Value = 28.83006510501794
$ws.Range("J7").Value = 10.31242315700045
$ws.Range("L7").Value = 11.9338037893232
$ws.Range("M7").Value = 17.27138995496635
$ws.Range("B8").Value = 17.94501528444402
$ws.Range("C8").Value = 8.299622612932476
$ws.Range("D8").Value = 7.560918133500433
$ws.Range("F8").Value = 39.22574678943732
$ws.Range("G8").Value = 45.75057949435803
$ws.Range("H8").Value = 18.5189528413457
$ws.Range("I8").Value = 28.65685125443197
$ws.Range("J8").Value = 10.29057254488222
$ws.Range("L8").Value = 11.92348994836079
$ws.Range("M8").Value = 17.37711399264798
$ws.Range("B9").Value = 18.97807487669221
$ws.Range("C9").Value = 9.280052971100057
$ws.Range("D9").Value = 7.563125389130161
$ws.Range("F9").Value = 39.04861130752529
$ws.Range("G9").Value = 45.68368950865604
$ws.Range("H9").Value = 18.38932852328471
$ws.Range("I9").Value = 28.36759075832932
$ws.Range("J9").Value = 10.25193121158746
$ws.Range("L9").Value = 11.91417901892353
$ws.Range("M9").Value = 17.60140496868804
$ws.Range("B10").Value = 19.71390572726199
$ws.Range("C10").Value = 9.932103934200732
$ws.Range("D10").Value = 7.566034987569202
$ws.Range("F10").Value = 38.96858840386424
$ws.Range("G10").Value = 45.71019801487279
$ws.Range("H10").Value = 18.31261362339652
$ws.Range("I10").Value = 28.18601094917585
$ws.Range("J10").Value = 10.22607960387522
$ws.Range("L10").Value = 11.91394775001516
$ws.Range("M10").Value = 17.77561550969306
$ws.Range("B11").Value = 20.04210021936856
$ws.Range("C11").Value = 10.21354708500522
$ws.Range("D11").Value = 7.567635350421521
$ws.Range("F11").Value = 38.94315025691058
$ws.Range("G11").Value = 45.73877972567181
$ws.Range("H11").Value = 18.28176006046956
$ws.Range("I11").Value = 28.11016933367537
$ws.Range("J11").Value = 10.21486445148456
$ws.Range("L11").Value = 11.9152671204699
$ws.Range("M11").Value = 17.85670799176268
$ws.Range("B12").Value = 20.16532332174136
$ws.Range("C12").Value = 10.31791735502302
$ws.Range("D12").Value = 7.56828090049236
$ws.Range("F12").Value = 38.93509959816603
$ws.Range("G12").Value = 45.75198195414164
$ws.Range("H12").Value = 18.27066013422768
$ws.Range("I12").Value = 28.08242609001982
$ws.Range("J12").Value = 10.21069548504503
$ws.Range("L12").Value = 11.91597055391902
$ws.Range("M12").Value = 17.8876624914454
$ws.Range("B13").Value = 20.13883379816205
$ws.Range("C13").Value = 10.29553778677852
$ws.Range("D13").Value = 7.56814011609795
$ws.Range("F13").Value = 38.93676300414286
$ws.Range("G13").Value = 45.74903279697145
$ws.Range("H13").Value = 18.27302471198169
$ws.Range("I13").Value = 28.08835760323259
$ws.Range("J13").Value = 10.2115898849737
$ws.Range("L13").Value = 11.91581000935542
$ws.Range("M13").Value = 17.88098521989203
$ws.Range("B14").Value = 20.05225967979787
$ws.Range("C14").Value = 10.222178020078
$ws.Range("D14").Value = 7.56768766959045
$ws.Range("F14").Value = 38.94245618612297
$ws.Range("G14").Value = 45.73981819646368
$ws.Range("H14").Value = 18.28083515412975
$ws.Range("I14").Value = 28.10786728871257
$ws.Range("J14").Value = 10.21451990769162
$ws.Range("L14").Value = 11.91532091374866
$ws.Range("M14").Value = 17.85924980344404
$ws.Range("B15").Value = 19.99908956253014
$ws.Range("C15").Value = 10.17695513337714
$ws.Range("D15").Value = 7.56741567257075
$ws.Range("F15").Value = 38.94614962349371
$ws.Range("G15").Value = 45.73448383639327
$ws.Range("H15").Value = 18.2856953436488
$ws.Range("I15").Value = 28.11994480831497
$ws.Range("J15").Value = 10.21632477365497
$ws.Range("L15").Value = 11.91504783960334
$ws.Range("M15").Value = 17.84596778395134
$ws.Range("B16").Value = 19.69231568474547
$ws.Range("C16").Value = 9.913403400705031
$ws.Range("D16").Value = 7.565935946236189
$ws.Range("F16").Value = 38.97047197555465
$ws.Range("G16").Value = 45.70866293722484
$ws.Range("H16").Value = 18.31471153879571
$ws.Range("I16").Value = 28.19110376777968
$ws.Range("J16").Value = 10.2268234718303
$ws.Range("L16").Value = 11.91389010515547
$ws.Range("M16").Value = 17.77035144124171
$ws.Range("B17").Value = 19.50235922218152
$ws.Range("C17").Value = 9.747816929583651
$ws.Range("D17").Value = 7.565098876954933
$ws.Range("F17").Value = 38.98820547440967
$ws.Range("G17").Value = 45.69705760580956
$ws.Range("H17").Value = 18.33354940347661
$ws.Range("I17").Value = 28.23649188877577
$ws.Range("J17").Value = 10.23340336424472
$ws.Range("L17").Value = 11.9135440249997
$ws.Range("M17").Value = 17.72442230393587
$ws.Range("B18").Value = 19.39249225674577
$ws.Range("C18").Value = 9.651148357691316
$ws.Range("D18").Value = 7.564643484355623
$ws.Range("F18").Value = 38.99943708801568
$ws.Range("G18").Value = 45.69193774129625
$ws.Range("H18").Value = 18.34476500705215
$ws.Range("I18").Value = 28.26323395705616
$ws.Range("J18").Value = 10.23723925315056
$ws.Range("L18").Value = 11.91347912287888
$ws.Range("M18").Value = 17.69817970632444
$ws.Range("B19").Value = 19.35519229253637
$ws.Range("C19").Value = 9.618173655220257
$ws.Range("D19").Value = 7.564493781614173
$ws.Range("F19").Value = 39.0034169449421
$ws.Range("G19").Value = 45.69047121403848
$ws.Range("H19").Value = 18.34862771549535
$ws.Range("I19").Value = 28.27239747009333
$ws.Range("J19").Value = 10.23854684331266
$ws.Range("L19").Value = 11.91348021967999
$ws.Range("M19").Value = 17.68932496427648
$ws.Range("B20").Value = 19.52264431866636
$ws.Range("C20").Value = 9.765591832694367
$ws.Range("D20").Value = 7.565185288496895
$ws.Range("F20").Value = 38.9862108855652
$ws.Range("G20").Value = 45.69813202318736
$ws.Range("H20").Value = 18.33150467863901
$ws.Range("I20").Value = 28.2315943891734
$ws.Range("J20").Value = 10.23269761598145
$ws.Range("L20").Value = 11.91356698867425
$ws.Range("M20").Value = 17.72929359870236
$ws.Range("B21").Value = 20.07771817087104
$ws.Range("C21").Value = 10.24378559689209
$ws.Range("D21").Value = 7.567819493513762
$ws.Range("F21").Value = 38.94074098033621
$ws.Range("G21").Value = 45.74246017456361
$ws.Range("H21").Value = 18.27852518182733
$ws.Range("I21").Value = 28.10211028975069
$ws.Range("J21").Value = 10.21365717599845
$ws.Range("L21").Value = 11.91545905005156
$ws.Range("M21").Value = 17.8656274805245
$ws.Range("B22").Value = 20.43428195371119
$ws.Range("C22").Value = 10.54344647737198
$ws.Range("D22").Value = 7.569771367162661
$ws.Range("F22").Value = 38.92024803414591
$ws.Range("G22").Value = 45.78529686370704
$ws.Range("H22").Value = 18.24730286923567
$ws.Range("I22").Value = 28.02317774285186
$ws.Range("J22").Value = 10.20166743122746
$ws.Range("L22").Value = 11.9178830920141
$ws.Range("M22").Value = 17.95615773336819
$ws.Range("B23").Value = 20.24458156205421
$ws.Range("C23").Value = 10.38469552711177
$ws.Range("D23").Value = 7.568708636216535
$ws.Range("F23").Value = 38.93033990249053
$ws.Range("G23").Value = 45.76116517786006
$ws.Range("H23").Value = 18.26365478245104
$ws.Range("I23").Value = 28.06478322244114
$ws.Range("J23").Value = 10.20802514269252
$ws.Range("L23").Value = 11.91648104114444
$ws.Range("M23").Value = 17.90771556078315
$ws.Range("B24").Value = 19.51347546852778
$ws.Range("C24").Value = 9.757560378186279
$ws.Range("D24").Value = 7.56514614132228
$ws.Range("F24").Value = 38.98710941121958
$ws.Range("G24").Value = 45.69764144424447
$ws.Range("H24").Value = 18.33242789813332
$ws.Range("I24").Value = 28.23380653100818
$ws.Range("J24").Value = 10.23301651958102
$ws.Range("L24").Value = 11.91355618912263
$ws.Range("M24").Value = 17.7270907781973
$ws.Range("B25").Value = 18.70212510385227
$ws.Range("C25").Value = 9.026660735606756
$ws.Range("D25").Value = 7.562301077690593
$ws.Range("F25").Value = 39.08776042634338
$ws.Range("G25").Value = 45.68854022448301
$ws.Range("H25").Value = 18.4211514549391
$ws.Range("I25").Value = 28.44042520071435
$ws.Range("J25").Value = 10.26193705213122
$ws.Range("L25").Value = 11.91553372237188
$ws.Range("M25").Value = 17.53900458636246
